$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRIP700A.xpc")

# New row 16 values, mirroring the structure of row 15 (HKL index 13 -> "HexGrid-60degTilt5degRes")
$ws.Range("A16").Value = 14
$ws.Range("A16").Style = $ws.Range("A15").Style

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$values = @(
    1.175415591871771,
    1.160408762032967,
    0.9646712100274477,
    0.9831432012714794,
    1.175415591871771,
    1.160408762032967,
    0.9601667918494853,
    0.9213191143260228,
    1.038141564909271,
    0.9832396425783597,
    1.174892722463677,
    1.062539986030207,
    1.070909691300916,
    1.02331323485835
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 3 + $i  # Column C = 3
    $ws.Cells.Item(16, $col).Value = $values[$i]
}
